# Updated symbol list on Mon Jan 16 11:43:45 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) values for the crypto
# listing on the active worksheet. Values are written with a leading
# apostrophe so Excel keeps them as literal text (matching the existing
# inline-string cells) instead of auto-converting numeric-looking strings
# into Number/Percentage values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'298.74"
$ws.Range("E2").Value = "'0.74%"
$ws.Range("D3").Value = "'31.26"
$ws.Range("E3").Value = "'0.03%"
$ws.Range("D4").Value = "'5.164"
$ws.Range("E4").Value = "'1.04%"
$ws.Range("D5").Value = "'0.08041"
$ws.Range("E5").Value = "'9.57%"
$ws.Range("D6").Value = "'2.673"
$ws.Range("E6").Value = "'60.19%"
$ws.Range("D7").Value = "'7.835"
$ws.Range("E7").Value = "'1.51%"
$ws.Range("E8").Value = "'2.51%"
$ws.Range("D9").Value = "'0.9182"
$ws.Range("E9").Value = "'-0.21%"
$ws.Range("D10").Value = "'0.1733"
$ws.Range("E10").Value = "'3.25%"
$ws.Range("D11").Value = "'0.07309"
$ws.Range("E11").Value = "'4.06%"
$ws.Range("D12").Value = "'0.08405"
$ws.Range("E12").Value = "'3.87%"
$ws.Range("D13").Value = "'0.03037"
$ws.Range("E13").Value = "'1.52%"
$ws.Range("D14").Value = "'0.09956"
$ws.Range("E14").Value = "'0.52%"
$ws.Range("D15").Value = "'0.001503"
$ws.Range("E15").Value = "'0.59%"
$ws.Range("D16").Value = "'0.005999"
$ws.Range("E16").Value = "'-2.70%"
$ws.Range("E17").Value = "'1.53%"
$ws.Range("D18").Value = "'2.246"
$ws.Range("E18").Value = "'0.84%"
$ws.Range("E19").Value = "'0.37%"
$ws.Range("E20").Value = "'-0.70%"
$ws.Range("D21").Value = "'4.629"
$ws.Range("E21").Value = "'1.60%"
$ws.Range("E22").Value = "'3.29%"
$ws.Range("D23").Value = "'0.04555"
$ws.Range("E23").Value = "'-2.03%"
$ws.Range("E24").Value = "'3.52%"
$ws.Range("E25").Value = "'0.42%"
$ws.Range("D26").Value = "'0.0001179"
$ws.Range("E26").Value = "'-9.14%"
$ws.Range("E27").Value = "'83.22%"
$ws.Range("D39").Value = "'0.01827"
$ws.Range("E39").Value = "'7.44%"
$ws.Range("D40").Value = "'0.04516"
$ws.Range("E40").Value = "'2.10%"
$ws.Range("D41").Value = "'0.007018"
$ws.Range("E41").Value = "'-2.14%"
$ws.Range("D42").Value = "'0.1342"
$ws.Range("E42").Value = "'0.94%"
$ws.Range("D43").Value = "'0.002239"
$ws.Range("E43").Value = "'4.77%"
$ws.Range("D44").Value = "'0.009814"
$ws.Range("E44").Value = "'-11.61%"
$ws.Range("D45").Value = "'0.00006463"
$ws.Range("E45").Value = "'8.03%"
$ws.Range("E46").Value = "'-0.07%"
$ws.Range("E47").Value = "'-39.28%"
$ws.Range("E48").Value = "'-56.68%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.07%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'0.00%"
